$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 275, shifting existing rows 275-361 down to 276-362
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new data record
$ws.Range("A275").Value = 11
$ws.Range("B275").Value = "Vega Monumental Concepción"
$ws.Range("C275").Value = "Bíobío"
$ws.Range("D275").Value = 44483
$ws.Range("E275").Value = 8
$ws.Range("F275").Value = 100114001
$ws.Range("G275").Value = "Papa"
$ws.Range("H275").Value = "Asterix"
$ws.Range("I275").Value = "1a (cosecha lavada)"
$ws.Range("J275").Value = 900
$ws.Range("K275").Value = 12000
$ws.Range("L275").Value = 14000
$ws.Range("M275").Value = 12889
$ws.Range("N275").Value = '$/malla 25 kilos'
$ws.Range("O275").Value = "Provincia de Melipilla"
$ws.Range("P275").Value = 516
$ws.Range("Q275").Value = 25
$ws.Range("R275").Value = "Hortaliza"
